# Insert a new data row at row 345 (pushing the existing rows 345-449 down
# to 346-450), then populate the new row 345 with its own data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A345").EntireRow.Insert()

$ws.Range("A345").Value = 10
$ws.Range("B345").Value = "Vega Modelo de Temuco"
$ws.Range("C345").Value = "La Araucanía"
$ws.Range("D345").Value = 44985
$ws.Range("E345").Value = 9
$ws.Range("F345").Value = 100112044
$ws.Range("G345").Value = "Perejil"
$ws.Range("H345").Value = "Sin especificar"
$ws.Range("I345").Value = "Primera"
$ws.Range("J345").Value = 50
$ws.Range("K345").Value = 5000
$ws.Range("L345").Value = 5000
$ws.Range("M345").Value = 5000
$ws.Range("N345").Value = "$/docena de atados (3 kilos)"
$ws.Range("O345").Value = "Provincia de Cautín"
$ws.Range("P345").Value = 1667
$ws.Range("Q345").Value = 3
$ws.Range("R345").Value = "Hortaliza"
